$wb = $excel.ActiveWorkbook

# 1) Add the new worksheet "Livres théorie mathématique" at the end of the workbook
$newSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "Livres théorie mathématique"
$newSheet.Range("A1").Value = "Livre"

# 2) Update the COURS sheet
$ws = $wb.Worksheets.Item("COURS")

# Renumber the existing chapter headers
$ws.Range("A12").Value = 2
$ws.Range("A23").Value = 3

# Copy the block format (rows 12-21) down to rows 34-43 for the new chapter entry
$ws.Range("A12:D12").Copy()
$ws.Range("A34").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A13:E21").Copy()
$ws.Range("A35").PasteSpecial(-4122)
$ws.Range("A13:E21").Copy()
$ws.Range("A35").PasteSpecial(-4163)
$ws.Application.CutCopyMode = $false

$ws.Rows.Item(34).RowHeight = $ws.Rows.Item(12).RowHeight

$ws.Range("A34").Value = 4
$ws.Range("B34").Value = "Recherche opérationnelle et applications - Bernard Fortz`n"

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 24
$ws.Range("B34").Select() | Out-Null
